$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 166; this shifts the existing rows 166-190
# down to 167-191 (old row 190's data ends up duplicated into new row 191).
$ws.Rows.Item(166).Insert()

# Populate the newly inserted row 166 with the weekly record, following
# the same template as the surrounding rows but with its own values.
$ws.Cells.Item(166, 1).Value = 7
$ws.Cells.Item(166, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(166, 3).Value = "Ñuble"
$ws.Cells.Item(166, 4).Value = 44491
$ws.Cells.Item(166, 5).Value = 16
$ws.Cells.Item(166, 6).Value = 100114013
$ws.Cells.Item(166, 7).Value = "Zanahoria"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 120
$ws.Cells.Item(166, 11).Value = 8500
$ws.Cells.Item(166, 12).Value = 9000
$ws.Cells.Item(166, 13).Value = 8750
$ws.Cells.Item(166, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(166, 15).Value = "Región de Ñuble"
$ws.Cells.Item(166, 16).Value = 438
$ws.Cells.Item(166, 17).Value = 20
$ws.Cells.Item(166, 18).Value = "Hortaliza"
